# Populate the "NovemberRaw" hidden data sheet with the final November 2024 statistics.
# This mirrors how the other *Raw sheets (e.g. OctoberRaw) are laid out: column A holds
# the library name, B/C/D hold the three monthly circulation metrics. The visible
# "November" sheet pulls these values via formulas, and "Yearly total" sums across all
# months, so populating NovemberRaw cascades through the whole workbook on recalculation.
$wb = $excel.ActiveWorkbook
$wsRaw = $wb.Worksheets.Item("NovemberRaw")

# Header row
$wsRaw.Cells.Item(1,1).Value = "Library"
$wsRaw.Cells.Item(1,2).Value = "Items owned by this library checked out at this library this month"
$wsRaw.Cells.Item(1,3).Value = "Items owned by other libraries checked out at this library this month"
$wsRaw.Cells.Item(1,4).Value = "Total circulation this month"

$wsRaw.Cells.Item(2,1).Value = "Atchison Public Library"
$wsRaw.Cells.Item(2,2).Value = 4099
$wsRaw.Cells.Item(2,3).Value = 1496
$wsRaw.Cells.Item(2,4).Value = 5595

$wsRaw.Cells.Item(3,1).Value = "Baldwin City Public Library"
$wsRaw.Cells.Item(3,2).Value = 2423
$wsRaw.Cells.Item(3,3).Value = 501
$wsRaw.Cells.Item(3,4).Value = 2924

$wsRaw.Cells.Item(4,1).Value = "Basehor Community Library"
$wsRaw.Cells.Item(4,2).Value = 6804
$wsRaw.Cells.Item(4,3).Value = 1088
$wsRaw.Cells.Item(4,4).Value = 7892

$wsRaw.Cells.Item(5,1).Value = "Bern Community Library"
$wsRaw.Cells.Item(5,2).Value = 87
$wsRaw.Cells.Item(5,3).Value = 40
$wsRaw.Cells.Item(5,4).Value = 127

$wsRaw.Cells.Item(6,1).Value = "Bonner Springs City Library"
$wsRaw.Cells.Item(6,2).Value = 4170
$wsRaw.Cells.Item(6,3).Value = 899
$wsRaw.Cells.Item(6,4).Value = 5069

$wsRaw.Cells.Item(7,1).Value = "Burlingame Community Library"
$wsRaw.Cells.Item(7,2).Value = 388
$wsRaw.Cells.Item(7,3).Value = 128
$wsRaw.Cells.Item(7,4).Value = 516

$wsRaw.Cells.Item(8,1).Value = "Carbondale City Library"
$wsRaw.Cells.Item(8,2).Value = 486
$wsRaw.Cells.Item(8,3).Value = 133
$wsRaw.Cells.Item(8,4).Value = 619

$wsRaw.Cells.Item(9,1).Value = "Centralia Community Library"
$wsRaw.Cells.Item(9,2).Value = 178
$wsRaw.Cells.Item(9,3).Value = 29
$wsRaw.Cells.Item(9,4).Value = 207

$wsRaw.Cells.Item(10,1).Value = "Corning City Library"
$wsRaw.Cells.Item(10,2).Value = 97
$wsRaw.Cells.Item(10,3).Value = 5
$wsRaw.Cells.Item(10,4).Value = 102

$wsRaw.Cells.Item(11,1).Value = "Digital Content"

$wsRaw.Cells.Item(12,1).Value = "Doniphan County Library - Elwood"
$wsRaw.Cells.Item(12,2).Value = 84
$wsRaw.Cells.Item(12,3).Value = 11
$wsRaw.Cells.Item(12,4).Value = 95

$wsRaw.Cells.Item(13,1).Value = "Doniphan County Library - Highland"
$wsRaw.Cells.Item(13,2).Value = 135
$wsRaw.Cells.Item(13,3).Value = 75
$wsRaw.Cells.Item(13,4).Value = 210

$wsRaw.Cells.Item(14,1).Value = "Doniphan County Library - Troy"
$wsRaw.Cells.Item(14,2).Value = 401
$wsRaw.Cells.Item(14,3).Value = 200
$wsRaw.Cells.Item(14,4).Value = 601

$wsRaw.Cells.Item(15,1).Value = "Doniphan County Library - Wathena"
$wsRaw.Cells.Item(15,2).Value = 221
$wsRaw.Cells.Item(15,3).Value = 48
$wsRaw.Cells.Item(15,4).Value = 269

$wsRaw.Cells.Item(16,1).Value = "Effingham Community Library"
$wsRaw.Cells.Item(16,2).Value = 217
$wsRaw.Cells.Item(16,3).Value = 37
$wsRaw.Cells.Item(16,4).Value = 254

$wsRaw.Cells.Item(17,1).Value = "Eudora Community Library"
$wsRaw.Cells.Item(17,2).Value = 1394
$wsRaw.Cells.Item(17,3).Value = 552
$wsRaw.Cells.Item(17,4).Value = 1946

$wsRaw.Cells.Item(18,1).Value = "Everest, Barnes Reading Room"
$wsRaw.Cells.Item(18,2).Value = 55
$wsRaw.Cells.Item(18,3).Value = 87
$wsRaw.Cells.Item(18,4).Value = 142

$wsRaw.Cells.Item(19,1).Value = "Hiawatha, Morrill Public Library"
$wsRaw.Cells.Item(19,2).Value = 1497
$wsRaw.Cells.Item(19,3).Value = 444
$wsRaw.Cells.Item(19,4).Value = 1941

$wsRaw.Cells.Item(20,1).Value = "Highland Community College"
$wsRaw.Cells.Item(20,2).Value = 29
$wsRaw.Cells.Item(20,3).Value = 37
$wsRaw.Cells.Item(20,4).Value = 66

$wsRaw.Cells.Item(21,1).Value = "Holton, Beck-Bookman Library"
$wsRaw.Cells.Item(21,2).Value = 1502
$wsRaw.Cells.Item(21,3).Value = 500
$wsRaw.Cells.Item(21,4).Value = 2002

$wsRaw.Cells.Item(22,1).Value = "Horton Public Library"
$wsRaw.Cells.Item(22,2).Value = 184
$wsRaw.Cells.Item(22,3).Value = 33
$wsRaw.Cells.Item(22,4).Value = 217

$wsRaw.Cells.Item(23,1).Value = "Lansing Community Library"
$wsRaw.Cells.Item(23,2).Value = 1764
$wsRaw.Cells.Item(23,3).Value = 475
$wsRaw.Cells.Item(23,4).Value = 2239

$wsRaw.Cells.Item(24,1).Value = "Leavenworth Public Library"
$wsRaw.Cells.Item(24,2).Value = 8075
$wsRaw.Cells.Item(24,3).Value = 1337
$wsRaw.Cells.Item(24,4).Value = 9412

$wsRaw.Cells.Item(25,1).Value = "Linwood Community Library"
$wsRaw.Cells.Item(25,2).Value = 470
$wsRaw.Cells.Item(25,3).Value = 166
$wsRaw.Cells.Item(25,4).Value = 636

$wsRaw.Cells.Item(26,1).Value = "Louisburg Library"

$wsRaw.Cells.Item(27,1).Value = "Lyndon Carnegie Library"
$wsRaw.Cells.Item(27,2).Value = 413
$wsRaw.Cells.Item(27,3).Value = 214
$wsRaw.Cells.Item(27,4).Value = 627

$wsRaw.Cells.Item(28,1).Value = "McLouth Public Library"
$wsRaw.Cells.Item(28,2).Value = 239
$wsRaw.Cells.Item(28,3).Value = 112
$wsRaw.Cells.Item(28,4).Value = 351

$wsRaw.Cells.Item(29,1).Value = "Meriden-Ozawkie Public Library"
$wsRaw.Cells.Item(29,2).Value = 1442
$wsRaw.Cells.Item(29,3).Value = 661
$wsRaw.Cells.Item(29,4).Value = 2103

$wsRaw.Cells.Item(30,1).Value = "Northeast Kansas Library System"
$wsRaw.Cells.Item(30,2).Value = 8
$wsRaw.Cells.Item(30,3).Value = 13
$wsRaw.Cells.Item(30,4).Value = 21

$wsRaw.Cells.Item(31,1).Value = "Nortonville Public Library"
$wsRaw.Cells.Item(31,2).Value = 191
$wsRaw.Cells.Item(31,3).Value = 63
$wsRaw.Cells.Item(31,4).Value = 254

$wsRaw.Cells.Item(32,1).Value = "Osage City Library"
$wsRaw.Cells.Item(32,2).Value = 1198
$wsRaw.Cells.Item(32,3).Value = 437
$wsRaw.Cells.Item(32,4).Value = 1635

$wsRaw.Cells.Item(33,1).Value = "Osawatomie Public Library"
$wsRaw.Cells.Item(33,2).Value = 712
$wsRaw.Cells.Item(33,3).Value = 254
$wsRaw.Cells.Item(33,4).Value = 966

$wsRaw.Cells.Item(34,1).Value = "Oskaloosa Public Library"
$wsRaw.Cells.Item(34,2).Value = 342
$wsRaw.Cells.Item(34,3).Value = 148
$wsRaw.Cells.Item(34,4).Value = 490

$wsRaw.Cells.Item(35,1).Value = "Ottawa Library"
$wsRaw.Cells.Item(35,2).Value = 5311
$wsRaw.Cells.Item(35,3).Value = 855
$wsRaw.Cells.Item(35,4).Value = 6166

$wsRaw.Cells.Item(36,1).Value = "Overbrook Public Library"
$wsRaw.Cells.Item(36,2).Value = 709
$wsRaw.Cells.Item(36,3).Value = 133
$wsRaw.Cells.Item(36,4).Value = 842

$wsRaw.Cells.Item(37,1).Value = "Paola Free Library"
$wsRaw.Cells.Item(37,2).Value = 2674
$wsRaw.Cells.Item(37,3).Value = 462
$wsRaw.Cells.Item(37,4).Value = 3136

$wsRaw.Cells.Item(38,1).Value = "Perry-Lecompton Community Library"
$wsRaw.Cells.Item(38,2).Value = 131
$wsRaw.Cells.Item(38,3).Value = 31
$wsRaw.Cells.Item(38,4).Value = 162

$wsRaw.Cells.Item(39,1).Value = "Pomona Community Library"
$wsRaw.Cells.Item(39,2).Value = 160
$wsRaw.Cells.Item(39,3).Value = 84
$wsRaw.Cells.Item(39,4).Value = 244

$wsRaw.Cells.Item(40,1).Value = "Prairie Hills Schools - Axtell Public School"
$wsRaw.Cells.Item(40,2).Value = 324
$wsRaw.Cells.Item(40,3).Value = 18
$wsRaw.Cells.Item(40,4).Value = 342

$wsRaw.Cells.Item(41,1).Value = "Prairie Hills Schools - Sabetha Elementary School"
$wsRaw.Cells.Item(41,2).Value = 2230
$wsRaw.Cells.Item(41,3).Value = 109
$wsRaw.Cells.Item(41,4).Value = 2339

$wsRaw.Cells.Item(42,1).Value = "Prairie Hills Schools - Sabetha High School"
$wsRaw.Cells.Item(42,2).Value = 8
$wsRaw.Cells.Item(42,3).Value = 16
$wsRaw.Cells.Item(42,4).Value = 24

$wsRaw.Cells.Item(43,1).Value = "Prairie Hills Schools - Sabetha Middle School"
$wsRaw.Cells.Item(43,2).Value = 147
$wsRaw.Cells.Item(43,3).Value = 11
$wsRaw.Cells.Item(43,4).Value = 158

$wsRaw.Cells.Item(44,1).Value = "Prairie Hills Schools - Wetmore Academic Center (Permanently closed)"

$wsRaw.Cells.Item(45,1).Value = "Richmond Public Library"
$wsRaw.Cells.Item(45,2).Value = 358
$wsRaw.Cells.Item(45,3).Value = 52
$wsRaw.Cells.Item(45,4).Value = 410

$wsRaw.Cells.Item(46,1).Value = "Rossville Community Library"
$wsRaw.Cells.Item(46,2).Value = 1306
$wsRaw.Cells.Item(46,3).Value = 417
$wsRaw.Cells.Item(46,4).Value = 1723

$wsRaw.Cells.Item(47,1).Value = "Sabetha, Mary Cotton Library"
$wsRaw.Cells.Item(47,2).Value = 2199
$wsRaw.Cells.Item(47,3).Value = 857
$wsRaw.Cells.Item(47,4).Value = 3056

$wsRaw.Cells.Item(48,1).Value = "Seneca Free Library"
$wsRaw.Cells.Item(48,2).Value = 1245
$wsRaw.Cells.Item(48,3).Value = 172
$wsRaw.Cells.Item(48,4).Value = 1417

$wsRaw.Cells.Item(49,1).Value = "Silver Lake Library"
$wsRaw.Cells.Item(49,2).Value = 1106
$wsRaw.Cells.Item(49,3).Value = 660
$wsRaw.Cells.Item(49,4).Value = 1766

$wsRaw.Cells.Item(50,1).Value = "Tonganoxie Public Library"
$wsRaw.Cells.Item(50,2).Value = 2783
$wsRaw.Cells.Item(50,3).Value = 690
$wsRaw.Cells.Item(50,4).Value = 3473

$wsRaw.Cells.Item(51,1).Value = "Valley Falls, Delaware Township Library"
$wsRaw.Cells.Item(51,2).Value = 462
$wsRaw.Cells.Item(51,3).Value = 191
$wsRaw.Cells.Item(51,4).Value = 653

$wsRaw.Cells.Item(52,1).Value = "Wellsville City Library"
$wsRaw.Cells.Item(52,2).Value = 839
$wsRaw.Cells.Item(52,3).Value = 312
$wsRaw.Cells.Item(52,4).Value = 1151

$wsRaw.Cells.Item(53,1).Value = "Wetmore Public Library"
$wsRaw.Cells.Item(53,2).Value = 98
$wsRaw.Cells.Item(53,3).Value = 233
$wsRaw.Cells.Item(53,4).Value = 331

$wsRaw.Cells.Item(54,1).Value = "Williamsburg Community Library"
$wsRaw.Cells.Item(54,2).Value = 162
$wsRaw.Cells.Item(54,3).Value = 38
$wsRaw.Cells.Item(54,4).Value = 200

$wsRaw.Cells.Item(55,1).Value = "Winchester Public Library"
$wsRaw.Cells.Item(55,2).Value = 256
$wsRaw.Cells.Item(55,3).Value = 422
$wsRaw.Cells.Item(55,4).Value = 678

# Reflect the author's last on-screen selection: they had clicked cell C24 on the
# visible "November" sheet while reviewing the update, then returned focus to the
# "Yearly total" sheet (which stays the active tab) before saving.
$wsNovember = $wb.Worksheets.Item("November")
$wsNovember.Range("C24").Select()
$wsYearly = $wb.Worksheets.Item("Yearly total")
$wsYearly.Activate()
